$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = "No value given"
$ws.Range("D2").Value = "button 1`nbutton 2`n"
$ws.Range("E2").Value = "First Title"
$ws.Range("F2").Value = "Description 1"
$ws.Hyperlinks.Add($ws.Range("G2"), "https://images.pexels.com/photos/268533/pexels-photo-268533.jpeg?cs=srgb&dl=pexels-pixabay-268533.jpg&fm=jpg`n", [Type]::Missing, [Type]::Missing, "https://images.pexels.com/photos/268533/pexels-photo-268533.jpeg?cs=srgb&dl=pexels-pixabay-268533.jpg&fm=jpg`n")
